# The workbook reports NATMI ligand-receptor edge statistics between four
# clusters (ECs, FAPs, MuSCs, Resolving-Mac) for the Sema4b -> Dcbld2 pair,
# laid out as a 4x4 grid of rows (2..17): sending cluster changes every 4
# rows, target cluster cycles within each block of 4.
#
# The commit "update scripts wuth new tpm" re-ran the pipeline with a new
# TPM table, which changes the underlying per-cluster average/total
# ligand expression (columns G/H, keyed by sending cluster) and receptor
# expression (columns M/N, keyed by target cluster). Every other touched
# column (I/J/O/P specificities and Q/R/S/T edge weights/specificities)
# is a deterministic function of those per-cluster values, so we recompute
# them the same way the NATMI script does rather than hardcoding each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster order used by the 4x4 row grid (rows 2-17).
$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# New ligand (Sema4b) average expression value per sending cluster.
$ligandAvg = @{
    "ECs"           = 0.6738823333333332
    "FAPs"          = 2.720340333333333
    "MuSCs"         = 3.560888666666667
    "Resolving-Mac" = 2.916489666666667
}

# New receptor (Dcbld2) average expression value per target cluster.
$receptorAvg = @{
    "ECs"           = 8.947221000000001
    "FAPs"          = 34.70130533333333
    "MuSCs"         = 14.549723
    "Resolving-Mac" = 0.5826356666666667
}

# Ligand-expressing / receptor-expressing cell counts (columns E / K) are
# unchanged by this update (still 3 for every row), so total = average * 3.
$cellCount = 3

$ligandTotal = @{}
foreach ($c in $clusters) { $ligandTotal[$c] = $ligandAvg[$c] * $cellCount }

$receptorTotal = @{}
foreach ($c in $clusters) { $receptorTotal[$c] = $receptorAvg[$c] * $cellCount }

# Derived specificity = value / sum(value across all sending/target clusters).
$sumLigandAvg = 0; foreach ($c in $clusters) { $sumLigandAvg += $ligandAvg[$c] }
$sumLigandTotal = 0; foreach ($c in $clusters) { $sumLigandTotal += $ligandTotal[$c] }
$sumReceptorAvg = 0; foreach ($c in $clusters) { $sumReceptorAvg += $receptorAvg[$c] }
$sumReceptorTotal = 0; foreach ($c in $clusters) { $sumReceptorTotal += $receptorTotal[$c] }

$ligandAvgSpec = @{}
$ligandTotalSpec = @{}
$receptorAvgSpec = @{}
$receptorTotalSpec = @{}
foreach ($c in $clusters) {
    $ligandAvgSpec[$c] = $ligandAvg[$c] / $sumLigandAvg
    $ligandTotalSpec[$c] = $ligandTotal[$c] / $sumLigandTotal
    $receptorAvgSpec[$c] = $receptorAvg[$c] / $sumReceptorAvg
    $receptorTotalSpec[$c] = $receptorTotal[$c] / $sumReceptorTotal
}

# Walk the 4x4 grid (rows 2-17) writing ligand columns (G,H,I,J), receptor
# columns (M,N,O,P) and the edge columns (Q,R,S,T) derived from them.
$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $g = $ligandAvg[$sending]
        $h = $ligandTotal[$sending]
        $i = $ligandAvgSpec[$sending]
        $j = $ligandTotalSpec[$sending]

        $m = $receptorAvg[$target]
        $n = $receptorTotal[$target]
        $o = $receptorAvgSpec[$target]
        $p = $receptorTotalSpec[$target]

        $ws.Range("G$row").Value = $g
        $ws.Range("H$row").Value = $h
        $ws.Range("I$row").Value = $i
        $ws.Range("J$row").Value = $j

        $ws.Range("M$row").Value = $m
        $ws.Range("N$row").Value = $n
        $ws.Range("O$row").Value = $o
        $ws.Range("P$row").Value = $p

        $ws.Range("Q$row").Value = $g * $m
        $ws.Range("R$row").Value = $h * $n
        $ws.Range("S$row").Value = $i * $o
        $ws.Range("T$row").Value = $j * $p

        $row += 1
    }
}
